# Duplicate the last question row on the "hard" sheet, then make that
# sheet the active / selected one (mirrors what previously was true of
# the "medium" sheet).

$wb = $excel.ActiveWorkbook

$hard = $wb.Worksheets.Item("hard")

# Activate the "hard" sheet (it becomes the selected / active tab; this
# also clears the previously-selected "medium" tab automatically).
$hard.Activate()

# Copy row 4 values into the new row 5 on the "hard" sheet.
$hard.Range("A5").Value = $hard.Range("A4").Text
$hard.Range("B5").Value = $hard.Range("B4").Text
$hard.Range("C5").Value = $hard.Range("C4").Text
$hard.Range("D5").Value = $hard.Range("D4").Text
$hard.Range("E5").Value = $hard.Range("E4").Text
$hard.Range("F5").Value = 3

# Update the selection on the hard sheet to match the new data range.
$hard.Range("A4:F5").Select()
